$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47/48: coin name & link swap (Monero <-> VeChain) ---
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

# --- Price (column D) updates: force text via leading apostrophe, then strip quote-prefix styling ---
$ws.Range("D2").Value = "'66.770.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.101.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'576.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'177.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'3.100.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Value = "'6.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Value = "'36.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "'3.621.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'66.806.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'3.104.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'16.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'481.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'7.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'83.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Value = "'10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("D32").Value = "'27.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Value = "'0.0₃0935"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'48.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Value = "'49.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Value = "'2.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'2.801.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'375.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'135.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.0344"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Value = "'25.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'2.23"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("E27").Value = "  -4.95%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").Value = "  +1.68%  "
